# Applies the changes described by the commit:
# - Add a "Calibrated" note next to the Source on the About sheet (B4)
# - Update CRtPaL-profits capacity-response calibration values (B2:B18)

$wb = $excel.ActiveWorkbook

# --- About sheet: add "Calibrated" label next to the Source row ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B4").Value = "Calibrated"

# --- CRtPaL-profits sheet: update calibrated response values ---
$wsProfits = $wb.Worksheets.Item("CRtPaL-profits")

$wsProfits.Range("B2").Value = 0.5    # hard coal
$wsProfits.Range("B3").Value = 0.5    # natural gas steam turbine
$wsProfits.Range("B4").Value = 0.5    # natural gas combined cycle
$wsProfits.Range("B5").Value = 0.02   # nuclear
$wsProfits.Range("B6").Value = 0.02   # hydro
$wsProfits.Range("B7").Value = 0.7    # onshore wind
$wsProfits.Range("B8").Value = 0.75   # solar PV
$wsProfits.Range("B9").Value = 0.5    # solar thermal
$wsProfits.Range("B10").Value = 0.5   # biomass
$wsProfits.Range("B11").Value = 0.02  # geothermal
$wsProfits.Range("B12").Value = 0.02  # petroleum
$wsProfits.Range("B13").Value = 0.15  # natural gas peaker
$wsProfits.Range("B14").Value = 0.5   # lignite
$wsProfits.Range("B15").Value = 0.5   # offshore wind
$wsProfits.Range("B16").Value = 0     # crude oil
$wsProfits.Range("B17").Value = 0     # heavy or residual fuel oil
$wsProfits.Range("B18").Value = 0     # municipal solid waste
